# Initial frequencies. Initial playing around with T Tests and ANOVAs.
#
# Under the "EBPAS Statistics" heading, right after the paragraph that
# introduces the EBPAS scoring (the "The first step is to compute EBPAS
# total score..." paragraph), add a new standalone paragraph with the
# author's aside about whether the divergence Likert scores should be
# reversed.

$d = $word.ActiveDocument

$anchorText = "The first step is to compute EBPAS total score, along with its subsets for scores of openness, divergence, and education."
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">__**One question I had: should the divergence Likert scores be reversed?__</w:t></w:r></w:p>'

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$anchorText*") {
        # Create a fresh empty paragraph right after the anchor paragraph,
        # then stamp it with the exact OOXML for the new text so the run
        # formatting matches a plain, unstyled paragraph.
        $p.Range.InsertParagraphAfter()
        $newPara = $p.Next()
        $newPara.Range.InsertXML($newParaXml) | Out-Null
        $found = $true
        break
    }
}

if (-not $found) {
    Write-Output "WARNING: anchor paragraph not found; document unchanged"
}
